# Arbeitszeit_Pichler.xlsx edit: add two new diary rows (36, 37) documenting
# work on torrent hashing methods ("worked on torrent hashing methods").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 36: 2019-08-15, 4 Stunden, Programmieren, "Logging, Bugfixing" ---
$ws.Range("E35").Copy()
$ws.Range("E36").PasteSpecial(-4122)
$ws.Range("E36").Value = 43692
$ws.Range("F36").Value = 4
$ws.Range("G36").Value = "Stunden"
$ws.Range("H36").Value = "Programmieren"
$ws.Range("I36").Value = "Logging, Bugfixing"

# --- Row 37: 2019-08-16, 4 Stunden, Programmieren, "Eruieren der Moeglichkeiten..." ---
$ws.Range("E35").Copy()
$ws.Range("E37").PasteSpecial(-4122)
$ws.Range("E37").Value = 43693
$ws.Range("F37").Value = 4
$ws.Range("G37").Value = "Stunden"
$ws.Range("H37").Value = "Programmieren"
$ws.Range("I37").Value = "Eruieren der Möglichkeiten der Hash-Umwandlung von Torrent-Dateien"

# Update the saved selection/active cell to match the author's last click.
$ws.Range("J40").Select()
